# Regenerate the "K" column (column G) on the active worksheet with
# freshly calculated strikeout values (replacing the previous Strike#
# derived figures), per the source commit:
#   "regen save_data to use K instead of Strike#, regen std/mean,
#    calc and write s_vals"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of worksheet row number -> new value for column G ("K")
$newK = @{
    2 = 2
    3 = 0
    4 = 0
    5 = 0
    6 = 0
    7 = 0
    8 = 1
    9 = 1
    10 = 1
    11 = 1
    12 = 0
    13 = 0
    14 = 1
    15 = 1
    16 = 2
    17 = 0
    19 = 1
    20 = 1
    21 = 1
    22 = 1
    23 = 0
    24 = 0
    25 = 2
    26 = 0
    27 = 2
    28 = 0
    29 = 0
    30 = 1
    31 = 3
    32 = 1
    33 = 0
    34 = 6
    35 = 0
    36 = 1
    37 = 0
    38 = 2
    39 = 1
    40 = 0
    41 = 3
    42 = 1
    43 = 0
    44 = 1
    45 = 2
    46 = 0
    47 = 1
    48 = 0
    49 = 1
    50 = 2
    51 = 1
    52 = 2
    53 = 2
    54 = 1
    56 = 3
    57 = 2
    60 = 1
}

# Column G is the 7th column ("K")
$col = 7

foreach ($row in $newK.Keys) {
    $ws.Cells.Item($row, $col).Value = $newK[$row]
}
